# Scheduled runner update: refresh market price / profit data across leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 556.5714
$ws.Range("I29").Value = 556.5714
$ws.Range("K29").Value = 1669.7142
$ws.Range("M29").Value = -1388.7142

$ws.Range("H115").Value = 4979
$ws.Range("I115").Value = 6965
$ws.Range("K115").Value = 20895
$ws.Range("M115").Value = -19328

$ws.Range("H128").Value = 38884.445
$ws.Range("J128").Value = 38884.445
$ws.Range("L128").Value = 38884.445
$ws.Range("N128").Value = -48844.445

$ws.Range("H130").Value = 50765
$ws.Range("J130").Value = 50765
$ws.Range("L130").Value = 50765
$ws.Range("N130").Value = -60805

$ws.Range("H132").Value = 246939.47
$ws.Range("I132").Value = 281092.3
$ws.Range("J132").Value = 1039
$ws.Range("K132").Value = 843276.8999999999
$ws.Range("L132").Value = 3117
$ws.Range("M132").Value = -840746.8999999999
$ws.Range("N132").Value = -8177

$ws.Range("H135").Value = 1843.1724
$ws.Range("I135").Value = 702.1667
$ws.Range("J135").Value = 7320
$ws.Range("K135").Value = 6319.5003
$ws.Range("L135").Value = 65880
$ws.Range("M135").Value = -3784.5003
$ws.Range("N135").Value = -70950

$ws.Range("H138").Value = 5141.6035
$ws.Range("I138").Value = 3096.0454
$ws.Range("J138").Value = 6391.6665
$ws.Range("K138").Value = 9288.136200000001
$ws.Range("L138").Value = 19174.9995
$ws.Range("M138").Value = -4148.136200000001
$ws.Range("N138").Value = -29454.9995

$ws.Range("H139").Value = 80780
$ws.Range("J139").Value = 80780
$ws.Range("L139").Value = 80780
$ws.Range("N139").Value = -91060

$ws.Range("H140").Value = 63779.5
$ws.Range("J140").Value = 63779.5
$ws.Range("L140").Value = 63779.5
$ws.Range("N140").Value = -74139.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6855.1113
$ws.Range("I74").Value = 726.13336
$ws.Range("K74").Value = 726.13336
$ws.Range("M74").Value = 147.86664

$ws.Range("H77").Value = 6855.1113
$ws.Range("I77").Value = 726.13336
$ws.Range("K77").Value = 3630.6668
$ws.Range("M77").Value = 737.3332

$ws.Range("H132").Value = 2030978.5
$ws.Range("I132").Value = 2465795.2
$ws.Range("J132").Value = 1833.3334
$ws.Range("K132").Value = 7397385.600000001
$ws.Range("L132").Value = 5500.0002
$ws.Range("M132").Value = -7394855.600000001
$ws.Range("N132").Value = -10560.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H134").Value = 45682.04
$ws.Range("I134").Value = 49284.824
$ws.Range("J134").Value = 4250
$ws.Range("K134").Value = 147854.472
$ws.Range("L134").Value = 12750
$ws.Range("M134").Value = -145319.472
$ws.Range("N134").Value = -17820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 264.23077
$ws.Range("I22").Value = 214.31429
$ws.Range("K22").Value = 214.31429
$ws.Range("M22").Value = 135.68571

$ws.Range("H31").Value = 1167.9269
$ws.Range("I31").Value = 961.4146
$ws.Range("J31").Value = 1374.439
$ws.Range("K31").Value = 961.4146
$ws.Range("L31").Value = 1374.439
$ws.Range("M31").Value = -666.4146
$ws.Range("N31").Value = -1964.439

$ws.Range("H34").Value = 1167.9269
$ws.Range("I34").Value = 961.4146
$ws.Range("J34").Value = 1374.439
$ws.Range("K34").Value = 961.4146
$ws.Range("L34").Value = 1374.439
$ws.Range("M34").Value = -759.4146
$ws.Range("N34").Value = -1778.439

$ws.Range("H132").Value = 10502.5
$ws.Range("I132").Value = 25012
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 75036
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -72506
$ws.Range("N132").Value = -22058

$ws.Range("H134").Value = 2465.8823
$ws.Range("I134").Value = 2512.2964
$ws.Range("J134").Value = 2286.8572
$ws.Range("K134").Value = 7536.889200000001
$ws.Range("L134").Value = 6860.571599999999
$ws.Range("M134").Value = -5001.889200000001
$ws.Range("N134").Value = -11930.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1274381.8
$ws.Range("J64").Value = 1751899.9
$ws.Range("L64").Value = 5255699.699999999
$ws.Range("N64").Value = -5256239.699999999

$ws.Range("H67").Value = 1274381.8
$ws.Range("J67").Value = 1751899.9
$ws.Range("L67").Value = 5255699.699999999
$ws.Range("N67").Value = -5257571.699999999

$ws.Range("H68").Value = 924.7347
$ws.Range("I68").Value = 679.5282999999999
$ws.Range("J68").Value = 1213.5333
$ws.Range("K68").Value = 2038.5849
$ws.Range("L68").Value = 3640.5999
$ws.Range("M68").Value = -1227.5849
$ws.Range("N68").Value = -5262.5999

$ws.Range("H71").Value = 924.7347
$ws.Range("I71").Value = 679.5282999999999
$ws.Range("J71").Value = 1213.5333
$ws.Range("K71").Value = 6115.7547
$ws.Range("L71").Value = 10921.7997
$ws.Range("M71").Value = -2059.7547
$ws.Range("N71").Value = -19033.7997

$ws.Range("H100").Value = 2660
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 9000
$ws.Range("N100").Value = -10622

$ws.Range("H105").Value = 227002260
$ws.Range("J105").Value = 227002260
$ws.Range("L105").Value = 681006780
$ws.Range("N105").Value = -681012022

$ws.Range("H107").Value = 56237.527
$ws.Range("J107").Value = 200840.6
$ws.Range("L107").Value = 602521.8
$ws.Range("N107").Value = -606361.8

$ws.Range("H131").Value = 1372354.6
$ws.Range("J131").Value = 1641257.4
$ws.Range("L131").Value = 4923772.199999999
$ws.Range("N131").Value = -4933852.199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4796.9653
$ws.Range("I70").Value = 4815.7
$ws.Range("J70").Value = 4755.3335
$ws.Range("K70").Value = 4815.7
$ws.Range("L70").Value = 4755.3335
$ws.Range("M70").Value = -4545.7
$ws.Range("N70").Value = -5295.3335

$ws.Range("H73").Value = 4796.9653
$ws.Range("I73").Value = 4815.7
$ws.Range("J73").Value = 4755.3335
$ws.Range("K73").Value = 4815.7
$ws.Range("L73").Value = 4755.3335
$ws.Range("M73").Value = -3879.7
$ws.Range("N73").Value = -6627.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8824.842000000001
$ws.Range("I132").Value = 14528.667
$ws.Range("J132").Value = 3691.4
$ws.Range("K132").Value = 43586.001
$ws.Range("L132").Value = 11074.2
$ws.Range("M132").Value = -41056.001
$ws.Range("N132").Value = -16134.2

$ws.Range("H136").Value = 1641.825
$ws.Range("I136").Value = 1444.2759
$ws.Range("J136").Value = 2162.6365
$ws.Range("K136").Value = 1444.2759
$ws.Range("L136").Value = 6487.9095
$ws.Range("M136").Value = -1782.8277
$ws.Range("N136").Value = -11587.9095
